# Updated cryptos list with GitHub Actions
# Refreshes the scraped Price (col D) and Volume(1h) (col E) figures
# for each coin row to the latest snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.466.15'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.28%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.836.40'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.60%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '260.32'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  -1.44%  '
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +0.17%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.5320'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +1.91%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3014'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -6.52%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06867'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +0.94%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '17.66'
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -6.01%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.844.86'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +0.10%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.7339'
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -5.86%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07360'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -5.23%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '88.96'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.960'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '13.91'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -0.45%  '
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +0.16%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007894'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.85%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '26.475.73'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.36%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '4.580'
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  -1.07%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.957'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -0.83%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.242'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -2.38%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '142.67'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.32%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.215'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +1.92%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.684'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.08%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '16.92'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.57%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '110.26'
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  -1.32%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '4.238'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +1.24%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.08806'
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +0.78%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.022'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -2.29%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.04797'
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -0.84%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.7302'
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +1.43%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.919'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +2.00%  '
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  -0.09%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '3.091'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '2.283'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +2.89%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01709'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -4.58%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.4715'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -3.11%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.9056'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +0.94%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '107.61'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -3.00%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.881'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -2.41%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.18%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '7.367'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -3.57%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.989'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -0.94%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4081'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -3.00%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.1229'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -1.01%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '34.81'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -0.53%  '
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -1.59%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.8919'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '60.02'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +0.23%  '
